$p = $ppt.ActivePresentation
$fonts = $p.Fonts
Write-Host "Fonts.Count: $($fonts.Count)"
for ($i=1; $i -le $fonts.Count; $i++) {
  $f = $fonts.Item($i)
  Write-Host "font $i : $($f.NameOther) idx=$($f.Index) id=$($f.Id)"
}
